$d = $word.ActiveDocument

# 1. Heading: warning sign -> [!] (capital "Review")
$d.Content.Find.Execute('⚠ Needs Review', $true, $false, $false, $false, $false, $true, 1, $false, '[!] Needs Review', 2) | Out-Null

# 2. Body tag: warning sign -> [!] (lowercase "review"), wrapped in brackets
$d.Content.Find.Execute('[⚠ Needs review]', $true, $false, $false, $false, $false, $true, 1, $false, '[[!] Needs review]', 2) | Out-Null

# 3. Q4 answer: reword the SSO/adaptability sentence
$d.Content.Find.Execute('The proposed technology solution, APX Stream, is a web-based platform accessible through standard web browsers, which significantly simplifies deployment and maintenance. This design ensures that APX Stream is inherently adaptable to changes in the software environment, eliminating the typical delays associated with software compatibility updates. However, it is important to note that the current version of APX Stream does not support Single-Sign-On (SSO), which may be a consideration for environments requiring streamlined user authentication processes. We are committed to continuously enhancing our platform and are considering the integration of SSO in future updates to meet our clients'' evolving needs.', $true, $false, $false, $false, $false, $true, 1, $false, 'The proposed technology solution, APX Stream, is a web-based platform accessible through standard web browsers, which significantly simplifies deployment and maintenance. This design ensures that APX Stream is inherently adaptable to various software environments without the extensive compatibility adjustments often required by traditional software installations. However, it is important to note that the current version of APX Stream does not support Single-Sign-On (SSO), which may be a consideration for environments requiring streamlined user authentication processes. We are committed to continuously enhancing our platform and are considering the integration of SSO in future updates to meet our clients'' evolving needs.', 2) | Out-Null

# 4. Q5 answer paragraph: rewrite each of the 8 content lines in place
$d.Content.Find.Execute('Our proposed technology solution, the APX Stream DataDrive, is designed to enhance operational efficiencies and deliver superior performance compared to in-house built systems and other market offerings. The key advantages of the APX Stream DataDrive include:', $true, $false, $false, $false, $false, $true, 1, $false, 'Our proposed technology solution, the APX Stream DataDrive, is designed to deliver superior operational efficiencies and enhanced performance compared to in-house built systems and other market offerings. Key advantages include:', 2) | Out-Null
$d.Content.Find.Execute('1. **Advanced Integration Capabilities**: Unlike many in-house systems that may struggle with integration complexities, our solution is built with state-of-the-art integration capabilities that ensure seamless connectivity with existing enterprise systems. This reduces the need for extensive customization and decreases integration costs and timelines.', $true, $false, $false, $false, $false, $true, 1, $false, '1. **Advanced Integration Capabilities**: The APX Stream DataDrive seamlessly integrates with existing IT infrastructure, allowing for a smoother and more efficient implementation process. This integration capability significantly reduces the risk of disruptions to ongoing operations, a common challenge with in-house systems which may not always align perfectly with other enterprise applications.', 2) | Out-Null
$d.Content.Find.Execute('2. **Real-Time Data Processing**: The APX Stream DataDrive excels in handling large volumes of data with minimal latency. This capability allows for real-time data analysis and decision-making, a critical component for dynamic and fast-paced business environments. This is a significant improvement over many in-house systems, which often process data in batches, leading to delays and outdated information.', $true, $false, $false, $false, $false, $true, 1, $false, '2. **Real-Time Data Processing**: Unlike many in-house solutions that may struggle with latency issues, our system ensures real-time data processing. This capability enables immediate insights and decision-making, crucial for maintaining competitive advantage in fast-paced industries.', 2) | Out-Null
$d.Content.Find.Execute('3. **Scalability and Flexibility**: Our technology is designed to grow with your business. It can easily scale up or down based on your needs without significant additional investment. This flexibility is often lacking in in-house systems, which might require substantial redevelopment to cope with changes in business scale or scope.', $true, $false, $false, $false, $false, $true, 1, $false, '3. **Continuous Improvement and Support**: Our solution comes with ongoing support and updates, ensuring that the system evolves in line with the latest technological advancements and security protocols. In contrast, in-house systems often require additional internal resources and time to update, which can divert focus from core business activities.', 2) | Out-Null
$d.Content.Find.Execute('4. **Enhanced Security Features**: APX Stream DataDrive employs robust security protocols that exceed industry standards. Our focus on security ensures that your data is protected against emerging threats, which is often a challenge for in-house systems to maintain over time without substantial investments in security updates.', $true, $false, $false, $false, $false, $true, 1, $false, '4. **Scalability and Flexibility**: The APX Stream DataDrive is designed to grow with your business, providing scalable solutions that adapt to increasing data volumes and complex operational demands without the need for significant additional investments.', 2) | Out-Null
$d.Content.Find.Execute('5. **Cost-Effectiveness**: By leveraging our solution, your firm can avoid the high upfront costs associated with developing, testing, and maintaining an in-house system. Additionally, our subscription model includes ongoing updates and support, ensuring that your system is always at the cutting edge without unexpected costs.', $true, $false, $false, $false, $false, $true, 1, $false, '5. **Enhanced Customer Satisfaction**: We employ open communications protocols, which facilitate immediate and ongoing customer feedback. This direct feedback loop allows for quicker resolutions and adaptations, enhancing customer satisfaction and loyalty.', 2) | Out-Null
$d.Content.Find.Execute('6. **Proven Customer Satisfaction**: We implement continuous feedback mechanisms through open communications protocols, allowing for immediate and ongoing customer satisfaction assessments. This approach helps in quickly identifying and addressing any issues, ensuring a high level of user satisfaction and system reliability.', $true, $false, $false, $false, $false, $true, 1, $false, '6. **Cost Efficiency**: By reducing the need for extensive in-house development and maintenance teams, our solution offers a cost-effective alternative to building and managing your own system. This translates into lower overall IT expenditures and a better allocation of resources towards strategic initiatives.', 2) | Out-Null
$d.Content.Find.Execute('In summary, the APX Stream DataDrive not only addresses the typical limitations of in-house built systems but also offers enhanced capabilities that drive operational efficiencies, reduce costs, and improve data security and handling. These advantages make our solution a compelling choice for firms looking to leverage technology for competitive advantage.', $true, $false, $false, $false, $false, $true, 1, $false, 'In summary, the APX Stream DataDrive not only addresses the typical limitations encountered with in-house built systems but also provides a robust framework for operational excellence and customer satisfaction. This comprehensive approach ensures that your investment in our technology yields measurable and sustainable benefits.', 2) | Out-Null

# 5. Remove Q6 through Q34 (question headers + answers) entirely,
#    leaving the document ending right after the Q5 answer paragraph.
$q6Para = $d.Paragraphs.Item(8)
$tailRange = $d.Range($q6Para.Range.Start, $d.Content.End)
$tailRange.Delete()

